# Update countries & provincias Spain
#
# The source data table (rows 206-214, column A = country name,
# column D = "Recuperados", column H = "Muertes") gets re-sorted: a few
# country names swap places with their neighbours, and the numeric
# values that go with each country move together with the country name.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column A: country names that change position --------------------
$ws.Range("A206").Value = "Islas Malvinas"
$ws.Range("A207").Value = "Groenlandia"
$ws.Range("A208").Value = "Islas Turcas y Caicos"
$ws.Range("A209").Value = "Santa Sede"
$ws.Range("A210").Value = "Seychelles"
$ws.Range("A211").Value = "Montserrat"
$ws.Range("A213").Value = "Papua Nueva Guinea"
$ws.Range("A214").Value = "Islas Virgenes Britanicas"

# --- Columns D (Recuperados) and H (Muertes) values that move with
#     the re-sorted countries -----------------------------------------
$ws.Range("D208").Value = 11
$ws.Range("H208").Value = 1

$ws.Range("D209").Value = 12
$ws.Range("H209").Value = 0

$ws.Range("D210").Value = 11
$ws.Range("H210").Value = 0

$ws.Range("D211").Value = 10
$ws.Range("H211").Value = 1

$ws.Range("D213").Value = 8
$ws.Range("H213").Value = 0

$ws.Range("D214").Value = 7
$ws.Range("H214").Value = 1
